# Auto-applied bulk data refresh across Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Updates cached market-price / profit columns (H:N) per scheduled runner pull.
$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "ALC" = @{
        "H41" = 290.3
        "I41" = 173.9
        "J41" = 406.7
        "K41" = 173.9
        "L41" = 406.7
        "M41" = 266.1
        "N41" = -1286.7
        "H86" = 1986.7778
        "I86" = 1940.1666
        "J86" = 2080
        "K86" = 1940.1666
        "L86" = 2080
        "M86" = -817.1666
        "N86" = -4326
        "H88" = 7833.4
        "I88" = 900
        "K88" = 900
        "M88" = -494
        "H89" = 1986.7778
        "I89" = 1940.1666
        "J89" = 2080
        "K89" = 9700.833000000001
        "L89" = 10400
        "M89" = -4084.833000000001
        "N89" = -21632
        "H91" = 7833.4
        "I91" = 900
        "K91" = 900
        "M91" = 504
        "H121" = 1922.5
        "J121" = 3000
        "L121" = 9000
        "N121" = -12494
        "H131" = 3380
        "I131" = 3380
        "K131" = 10140
        "M131" = -5100
        "H132" = 2433.5715
        "I132" = 2756.3635
        "J132" = 1250
        "K132" = 8269.0905
        "L132" = 3750
        "M132" = -5739.0905
        "N132" = -8810
        "H137" = 1607.8636
        "I137" = 1595.625
        "J137" = 1640.5
        "K137" = 4786.875
        "L137" = 4921.5
        "M137" = -2236.875
        "N137" = -10021.5
        "H141" = 2145.5833
        "I141" = 1471.8889
        "J141" = 4166.6665
        "K141" = 4415.6667
        "L141" = 12499.9995
        "M141" = 764.3333000000002
        "N141" = -22859.9995
    }
    "ARM" = @{
        "H34" = 9975
        "I34" = 8962.5
        "J34" = 12000
        "K34" = 8962.5
        "L34" = 12000
        "M34" = -8691.5
        "N34" = -12542
        "H61" = 2523.074
        "I61" = 2016.2632
        "J61" = 3726.75
        "K61" = 2016.2632
        "L61" = 3726.75
        "M61" = -1804.2632
        "N61" = -4150.75
        "H136" = 2523.074
        "I136" = 2016.2632
        "J136" = 3726.75
        "K136" = 6048.7896
        "L136" = 11180.25
        "M136" = -3498.7896
        "N136" = -16280.25
    }
    "BSM" = @{
        "H99" = 166667920
        "I99" = 333334430
        "J99" = 1396.6666
        "K99" = 333334430
        "L99" = 1396.6666
        "M99" = -333332932
        "N99" = -4392.6666
        "H107" = 814.375
        "I107" = 618.5454999999999
        "J107" = 1245.2
        "K107" = 618.5454999999999
        "L107" = 1245.2
        "M107" = 1301.4545
        "N107" = -5085.2
    }
    "CRP" = @{
        "H9" = 11710.667
        "J9" = 11710.667
        "L9" = 11710.667
        "N9" = -12046.667
        "H31" = 3458.5303
        "I31" = 1386.037
        "J31" = 12784.75
        "K31" = 1386.037
        "L31" = 12784.75
        "M31" = -1091.037
        "N31" = -13374.75
        "H34" = 3458.5303
        "I34" = 1386.037
        "J34" = 12784.75
        "K34" = 1386.037
        "L34" = 12784.75
        "M34" = -1184.037
        "N34" = -13188.75
        "H122" = 2781119
        "I122" = 5558535
        "J122" = 3702.8
        "K122" = 16675605
        "L122" = 11108.4
        "M122" = -16673155
        "N122" = -16008.4
        "H132" = 3217.3333
        "I132" = 3052.4285
        "J132" = 3361.625
        "K132" = 9157.2855
        "L132" = 10084.875
        "M132" = -6627.2855
        "N132" = -15144.875
        "H134" = 2619.4443
        "I134" = 2573.8718
        "K134" = 7721.6154
        "M134" = -5186.6154
    }
    "CUL" = @{
        "H114" = 5689.696
        "I114" = 510.8
        "K114" = 1532.4
        "M114" = 1721.6
        "H117" = 18525452
        "I117" = 14661.286
        "J117" = 30305044
        "K117" = 43983.858
        "L117" = 90915132
        "M117" = -40541.858
        "N117" = -90922016
        "H121" = 7869.6523
        "I121" = 370.64285
        "J121" = 19534.777
        "K121" = 1111.92855
        "L121" = 58604.33099999999
        "M121" = 198.0714499999999
        "N121" = -61224.33099999999
    }
    "GSM" = @{
        "H80" = 2971.158
        "I80" = 2740
        "J80" = 3289
        "K80" = 2740
        "L80" = 3289
        "M80" = -1742
        "N80" = -5285
        "H83" = 2971.158
        "I83" = 2740
        "J83" = 3289
        "K83" = 13700
        "L83" = 16445
        "M83" = -8708
        "N83" = -26429
    }
    "LTW" = @{
        "H40" = 40002324
        "I40" = 50002124
        "J40" = 3115
        "K40" = 50002124
        "L40" = 3115
        "M40" = -50001988
        "N40" = -3387
        "H46" = 611.34784
        "I46" = 464.33334
        "J46" = 705.8570999999999
        "K46" = 464.33334
        "L46" = 705.8570999999999
        "M46" = -276.33334
        "N46" = -1081.8571
        "H82" = 381949.78
        "I82" = 1001399.9
        "J82" = 55923.42
        "K82" = 1001399.9
        "L82" = 55923.42
        "M82" = -1001038.9
        "N82" = -56645.42
        "H85" = 381949.78
        "I85" = 1001399.9
        "J85" = 55923.42
        "K85" = 1001399.9
        "L85" = 55923.42
        "M85" = -1000151.9
        "N85" = -58419.42
        "H136" = 7098.207
        "I136" = 2931.6191
        "K136" = 8794.8573
        "M136" = -6244.8573
    }
    "WVR" = @{
        "H81" = 1652
        "I81" = 1173.5714
        "J81" = 2210.1667
        "K81" = 2347.1428
        "L81" = 4420.3334
        "M81" = -1286.1428
        "N81" = -6542.3334
        "H84" = 1652
        "I84" = 1173.5714
        "J84" = 2210.1667
        "K84" = 11735.714
        "L84" = 22101.667
        "M84" = -6431.714
        "N84" = -32709.667
        "H122" = 888.8484999999999
        "I122" = 924.5925999999999
        "J122" = 728
        "K122" = 2773.7778
        "L122" = 2184
        "M122" = -323.7777999999998
        "N122" = -7084
        "H126" = 1120.7778
        "I126" = 830.3333
        "J126" = 1701.6666
        "K126" = 2490.9999
        "L126" = 5104.9998
        "M126" = -20.9998999999998
        "N126" = -10044.9998
        "H132" = 1654.8572
        "I132" = 1317.3572
        "J132" = 1992.3572
        "K132" = 3952.0716
        "L132" = 5977.071599999999
        "M132" = -1422.0716
        "N132" = -11037.0716
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
